$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.15"
$ws.Range("D3").Value = "'22.99"
$ws.Range("D4").Value = "'5.381"
$ws.Range("D5").Value = "'0.05969"
$ws.Range("D6").Value = "'3.400"
$ws.Range("D7").Value = "'6.481"
$ws.Range("D8").Value = "'0.8084"
$ws.Range("D9").Value = "'0.9093"
$ws.Range("D11").Value = "'0.07411"
$ws.Range("D12").Value = "'0.03314"
$ws.Range("D13").Value = "'0.03054"
$ws.Range("D14").Value = "'0.09344"
$ws.Range("D15").Value = "'3.842"
$ws.Range("D16").Value = "'0.001571"
$ws.Range("D17").Value = "'0.04545"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006076"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005030"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "UpBots"
$ws.Range("C20").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D20").Value = "'0.007494"
$ws.Range("E20").Value = "19UpBotsUBXTBestin24h"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009841"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.00007787"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.617"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.163"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.01120"
$ws.Range("E25").Value = "24OneONE"
$ws.Range("D40").Value = "'0.03895"
$ws.Range("D41").Value = "'0.006153"
$ws.Range("D43").Value = "'0.002796"
$ws.Range("D44").Value = "'0.007165"
$ws.Range("D45").Value = "'0.00005191"
$ws.Range("D48").Value = "'1.043"
$ws.Range("D49").Value = "'0.002257"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("D51").Value = "'0.0001997"
